$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

# Column G holds "Recorded By". Some rows list the recorder as
# "System, dnasr281@gmail.com" - swap the order to "dnasr281@gmail.com, System".
$col = $ws.Range("G1:G319")
$col.Replace("System, dnasr281@gmail.com", "dnasr281@gmail.com, System") | Out-Null
